# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '30.550.85'
$ws.Range('E2').Value = '  -0.54%  '

# Row 3
$ws.Range('D3').Value = '1.877.94'

# Row 4
$ws.Range('D4').Value = '''0.9997'
$ws.Range('E4').Value = '  -0.34%  '

# Row 5
$ws.Range('D5').Value = '''236.28'
$ws.Range('E5').Value = '  -3.56%  '

# Row 6
$ws.Range('D6').Value = '''1.000'
$ws.Range('E6').Value = '  -0.24%  '

# Row 7
$ws.Range('D7').Value = '''0.4882'
$ws.Range('E7').Value = '  -1.45%  '

# Row 8
$ws.Range('D8').Value = '''0.2905'
$ws.Range('E8').Value = '  -1.80%  '

# Row 9
$ws.Range('D9').Value = '''0.06669'
$ws.Range('E9').Value = '  -1.95%  '

# Row 10
$ws.Range('D10').Value = '1.877.38'
$ws.Range('E10').Value = '  -0.66%  '

# Row 11
$ws.Range('D11').Value = '''16.62'
$ws.Range('E11').Value = '  -2.70%  '

# Row 12
$ws.Range('D12').Value = '''0.07236'
$ws.Range('E12').Value = '  -1.04%  '

# Row 13
$ws.Range('D13').Value = '''88.79'
$ws.Range('E13').Value = '  -2.30%  '

# Row 14
$ws.Range('D14').Value = '''5.001'
$ws.Range('E14').Value = '  -1.03%  '

# Row 15
$ws.Range('D15').Value = '''0.6516'
$ws.Range('E15').Value = '  -3.16%  '

# Row 16
$ws.Range('D16').Value = '30.509.51'
$ws.Range('E16').Value = '  -0.70%  '

# Row 17
$ws.Range('D17').Value = '''0.000007881'
$ws.Range('E17').Value = '  -1.22%  '

# Row 18
$ws.Range('D18').Value = '''1.000'
$ws.Range('E18').Value = '  -0.23%  '

# Row 19
$ws.Range('D19').Value = '''12.97'
$ws.Range('E19').Value = '  -2.22%  '

# Row 20
$ws.Range('D20').Value = '2.119.09'
$ws.Range('E20').Value = '  -1.23%  '

# Row 21
$ws.Range('D21').Value = '''1.000'
$ws.Range('E21').Value = '  -0.20%  '

# Row 22
$ws.Range('D22').Value = '''4.713'
$ws.Range('E22').Value = '  -2.84%  '

# Row 23
$ws.Range('D23').Value = '''193.55'
$ws.Range('E23').Value = '  +9.49%  '

# Row 24
$ws.Range('D24').Value = '''6.108'
$ws.Range('E24').Value = '  +0.68%  '

# Row 25
$ws.Range('D25').Value = '''9.316'
$ws.Range('E25').Value = '  +0.36%  '

# Row 26
$ws.Range('D26').Value = '''156.75'
$ws.Range('E26').Value = '  +0.90%  '

# Row 27
$ws.Range('E27').Value = '  -0.63%  '

# Row 28
$ws.Range('D28').Value = '''1.823'
$ws.Range('E28').Value = '  -5.61%  '

# Row 29
$ws.Range('D29').Value = '''1.404'
$ws.Range('E29').Value = '  +1.76%  '

# Row 30
$ws.Range('D30').Value = '''4.250'
$ws.Range('E30').Value = '  -1.96%  '

# Row 31
$ws.Range('D31').Value = '''0.09019'
$ws.Range('E31').Value = '  +1.25%  '

# Row 32
$ws.Range('D32').Value = '''3.922'
$ws.Range('E32').Value = '  -2.51%  '

# Row 33
$ws.Range('D33').Value = '''0.05124'
$ws.Range('E33').Value = '  -2.34%  '

# Row 34
$ws.Range('D34').Value = '''0.7233'
$ws.Range('E34').Value = '  -2.62%  '

# Row 35
$ws.Range('D35').Value = '''1.078'
$ws.Range('E35').Value = '  -4.95%  '

# Row 36
$ws.Range('D36').Value = '''2.693'
$ws.Range('E36').Value = '  +0.86%  '

# Row 37
$ws.Range('D37').Value = '''0.01812'
$ws.Range('E37').Value = '  -3.69%  '

# Row 38
$ws.Range('D38').Value = '''2.661'
$ws.Range('E38').Value = '  -1.30%  '

# Row 39
$ws.Range('D39').Value = '''0.9168'
$ws.Range('E39').Value = '  -2.16%  '

# Row 40
$ws.Range('D40').Value = '''2.046'
$ws.Range('E40').Value = '  -5.51%  '

# Row 41
$ws.Range('D41').Value = '''0.4390'
$ws.Range('E41').Value = '  +0.94%  '

# Row 42
$ws.Range('D42').Value = '''104.78'
$ws.Range('E42').Value = '  -0.41%  '

# Row 43
$ws.Range('D43').Value = '''0.9957'

# Row 44
$ws.Range('D44').Value = '''5.727'
$ws.Range('E44').Value = '  -1.44%  '

# Row 45
$ws.Range('D45').Value = '''0.1329'
$ws.Range('E45').Value = '  -2.21%  '

# Row 46
$ws.Range('D46').Value = '''7.381'
$ws.Range('E46').Value = '  -3.68%  '

# Row 47
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = '''0.4013'
$ws.Range('E47').Value = '  +3.40%  '

# Row 48
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '''0.05817'
$ws.Range('E48').Value = '  -0.35%  '

# Row 49
$ws.Range('D49').Value = '''8.671'
$ws.Range('E49').Value = '  +1.56%  '

# Row 50
$ws.Range('D50').Value = '''1.404'
$ws.Range('E50').Value = '  +1.69%  '

# Row 51
$ws.Range('D51').Value = '''33.12'
$ws.Range('E51').Value = '  -0.81%  '
